# Fix typo in geo-coordinate strings: commas used instead of decimal points
# in the latitude portion of two "Location (Geo)" values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "39.734928,-8.820685"
$ws.Range("B9").Value = "39.734905,-8.820718"
